$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "TP"
$ws.Range("B1").Value = "n"

$ws.Range("A2").Value = "T0"
$ws.Range("B2").Value = 8

$ws.Range("A3").Value = "T1"
$ws.Range("B3").Value = 17

$ws.Range("A4").Value = "T2"
$ws.Range("B4").Value = 17

$ws.Range("A5").Value = "T3"
$ws.Range("B5").Value = 19

$ws.Range("A6").Value = "T4"
$ws.Range("B6").Value = 11

$ws.Range("A7").Value = "T5"
$ws.Range("B7").Value = 14

$ws.Range("A8").Value = "TF"
$ws.Range("B8").Value = 7
